# Loan RBI, Variable Instalments
#
# - "Repayment schedule" sheet: insert a new (blank) column before column N,
#   shifting the old "Late" / "heading" / "Outstanding" columns one to the
#   right (N->O, O->P, P->Q). Give the new column a width matching the
#   neighbouring "In Advance" column.
# - Make "Repayment schedule" the active sheet/tab (instead of "Transactions"),
#   and update its selection.

$wb = $excel.ActiveWorkbook

$wsSchedule = $wb.Worksheets.Item("Repayment schedule")

# Insert a blank column in front of column N ("Late"), pushing the
# remaining columns (Late / heading / Outstanding) one column to the right.
[void]$wsSchedule.Columns("N").EntireColumn.Insert()

# Match the width of the newly inserted column to column M ("In Advance").
$wsSchedule.Columns("N").ColumnWidth = 10.7109375

# "Repayment schedule" becomes the active/selected sheet (this also clears
# the previously selected "Transactions" tab), with a new selection.
[void]$wsSchedule.Activate()
[void]$wsSchedule.Range("M17").Select()
